$wb = $excel.ActiveWorkbook

# --- OFF sheet: update divisional round row (row 3) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 284
$wsOff.Range("C3").Value = 199
$wsOff.Range("D3").Value = 66
$wsOff.Range("E3").Value = 38

# --- DEF sheet: update divisional round row (row 3) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 249
$wsDef.Range("C3").Value = 155
$wsDef.Range("D3").Value = 56
$wsDef.Range("E3").Value = 24

$wb.Save()
